$d = $word.ActiveDocument

# The document starts with six paragraphs:
#   1. "Eu,teste query teste query, ..." (explicit spacing-after, ends with a line break run)
#   2. empty (bare pPr)
#   3. empty (bottom-border pPr)
#   4. "Modelo Número 1" (explicit spacing-after, ends with a line break run)
#   5. empty (bare pPr)
#   6. empty (bottom-border pPr)
#
# The target keeps a single paragraph with an empty <w:pPr/> and one run
# reading "LIVRO N.° 2 - REGISTRO". Rather than stripping the formatting
# back off paragraph 1 (which would leave explicit zeroed overrides behind),
# drop paragraph 1 completely so the already-bare paragraph 2 becomes the
# sole survivor, then delete the remaining trailing paragraphs and fill in
# the new text.

# Remove paragraph 1 (its text, break run and paragraph mark) by deleting
# from the very start of the document through the end of its mark; this
# merges paragraph 2 up into position 1, carrying over its untouched,
# empty <w:pPr/>.
$d.Range(0, $d.Paragraphs(1).Range.End).Delete()

# Delete the remaining paragraphs (old paragraphs 3-6, now 2-5), working
# from the end forward. Each deleted range starts one character before the
# previous paragraph's end so the preceding paragraph mark is consumed too
# (a paragraph's own mark can't be removed by deleting only its own range).
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $prevEnd = $d.Paragraphs($i - 1).Range.End
    $thisEnd = $d.Paragraphs($i).Range.End
    $d.Range($prevEnd - 1, $thisEnd).Delete()
}

# Fill the surviving (still formatting-free) paragraph with the new text.
$d.Paragraphs(1).Range.Text = "LIVRO N.° 2 - REGISTRO"
